$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1591.1613
$ws.Range("I86").Value = 1591.1613
$ws.Range("J86").Value = 0
$ws.Range("K86").Value = 1591.1613
$ws.Range("L86").Value = 0
$ws.Range("M86").Value = -468.1613
$ws.Range("N86").ClearContents()
$ws.Range("H89").Value = 1591.1613
$ws.Range("I89").Value = 1591.1613
$ws.Range("J89").Value = 0
$ws.Range("K89").Value = 7955.8065
$ws.Range("L89").Value = 0
$ws.Range("M89").Value = -2339.8065
$ws.Range("N89").ClearContents()
$ws.Range("H106").Value = 73002
$ws.Range("I106").Value = 111670
$ws.Range("J106").Value = 15000
$ws.Range("K106").Value = 111670
$ws.Range("L106").Value = 15000
$ws.Range("M106").Value = -111039
$ws.Range("N106").Value = -16262
$ws.Range("H129").Value = 833.38464
$ws.Range("I129").Value = 630.6667
$ws.Range("J129").Value = 1007.1429
$ws.Range("K129").Value = 1892.0001
$ws.Range("L129").Value = 3021.4287
$ws.Range("M129").Value = 3107.9999
$ws.Range("N129").Value = -13021.4287
$ws.Range("H137").Value = 7265.5293
$ws.Range("I137").Value = 1055.25
$ws.Range("J137").Value = 12785.777
$ws.Range("K137").Value = 3165.75
$ws.Range("L137").Value = 38357.331
$ws.Range("M137").Value = -615.75
$ws.Range("N137").Value = -43457.331
$ws.Range("H138").Value = 114735.99
$ws.Range("I138").Value = 1144.5676
$ws.Range("J138").Value = 195560.66
$ws.Range("K138").Value = 3433.7028
$ws.Range("L138").Value = 586681.98
$ws.Range("M138").Value = 1706.2972
$ws.Range("N138").Value = -596961.98

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17031.283
$ws.Range("I32").Value = 11581.772
$ws.Range("J32").Value = 30864.654
$ws.Range("K32").Value = 11581.772
$ws.Range("L32").Value = 30864.654
$ws.Range("M32").Value = -11294.772
$ws.Range("N32").Value = -31438.654
$ws.Range("H45").Value = 1067.52
$ws.Range("I45").Value = 978.8570999999999
$ws.Range("J45").Value = 1180.3636
$ws.Range("K45").Value = 978.8570999999999
$ws.Range("L45").Value = 1180.3636
$ws.Range("M45").Value = -601.8570999999999
$ws.Range("N45").Value = -1934.3636
$ws.Range("H61").Value = 1699
$ws.Range("I61").Value = 1445.0769
$ws.Range("J61").Value = 5000
$ws.Range("K61").Value = 1445.0769
$ws.Range("L61").Value = 5000
$ws.Range("M61").Value = -1233.0769
$ws.Range("N61").Value = -5424
$ws.Range("H74").Value = 1514.1708
$ws.Range("I74").Value = 1383.9395
$ws.Range("K74").Value = 1383.9395
$ws.Range("M74").Value = -509.9395
$ws.Range("H77").Value = 1514.1708
$ws.Range("I77").Value = 1383.9395
$ws.Range("K77").Value = 6919.6975
$ws.Range("M77").Value = -2551.6975
$ws.Range("H136").Value = 1699
$ws.Range("I136").Value = 1445.0769
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 4335.2307
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -1785.2307
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1550.9445
$ws.Range("I20").Value = 1336.875
$ws.Range("J20").Value = 1722.2
$ws.Range("K20").Value = 1336.875
$ws.Range("L20").Value = 1722.2
$ws.Range("M20").Value = -1089.875
$ws.Range("N20").Value = -2216.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 18445.25
$ws.Range("I31").Value = 21018.117
$ws.Range("J31").Value = 13196.6
$ws.Range("K31").Value = 21018.117
$ws.Range("L31").Value = 13196.6
$ws.Range("M31").Value = -20723.117
$ws.Range("N31").Value = -13786.6
$ws.Range("H34").Value = 18445.25
$ws.Range("I34").Value = 21018.117
$ws.Range("J34").Value = 13196.6
$ws.Range("K34").Value = 21018.117
$ws.Range("L34").Value = 13196.6
$ws.Range("M34").Value = -20816.117
$ws.Range("N34").Value = -13600.6
$ws.Range("H62").Value = 3493.1667
$ws.Range("I62").Value = 2502.375
$ws.Range("K62").Value = 2502.375
$ws.Range("M62").Value = -1878.375
$ws.Range("H65").Value = 3493.1667
$ws.Range("I65").Value = 2502.375
$ws.Range("K65").Value = 12511.875
$ws.Range("M65").Value = -9391.875
$ws.Range("H132").Value = 26595.475
$ws.Range("I132").Value = 39807
$ws.Range("J132").Value = 2059.7856
$ws.Range("K132").Value = 119421
$ws.Range("L132").Value = 6179.3568
$ws.Range("M132").Value = -116891
$ws.Range("N132").Value = -11239.3568

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 6178.125
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 6178.125
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 18534.375
$ws.Range("M75").ClearContents()
$ws.Range("N75").Value = -20530.375
$ws.Range("H78").Value = 6178.125
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 6178.125
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 55603.125
$ws.Range("M78").ClearContents()
$ws.Range("N78").Value = -65587.125
$ws.Range("H93").Value = 2576.75
$ws.Range("J93").Value = 2576.75
$ws.Range("L93").Value = 7730.25
$ws.Range("N93").Value = -11474.25
$ws.Range("H113").Value = 764.6667
$ws.Range("I113").Value = 1000
$ws.Range("J113").Value = 647
$ws.Range("K113").Value = 3000
$ws.Range("L113").Value = 1941
$ws.Range("M113").Value = -830
$ws.Range("N113").Value = -6281
$ws.Range("H131").Value = 122797.05
$ws.Range("J131").Value = 145846.92
$ws.Range("L131").Value = 437540.76
$ws.Range("N131").Value = -447620.76
$ws.Range("H132").Value = 1856.5
$ws.Range("I132").Value = 1100
$ws.Range("J132").Value = 1951.0625
$ws.Range("K132").Value = 9900
$ws.Range("L132").Value = 17559.5625
$ws.Range("M132").Value = -7370
$ws.Range("N132").Value = -22619.5625

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 32083.705
$ws.Range("I40").Value = 2329.5715
$ws.Range("J40").Value = 39797.742
$ws.Range("K40").Value = 2329.5715
$ws.Range("L40").Value = 39797.742
$ws.Range("M40").Value = -2193.5715
$ws.Range("N40").Value = -40069.742

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 3165.932
$ws.Range("I132").Value = 589.5128
$ws.Range("J132").Value = 23262
$ws.Range("K132").Value = 1768.5384
$ws.Range("L132").Value = 69786
$ws.Range("M132").Value = 761.4616000000001
$ws.Range("N132").Value = -74846
